$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1961.9474
$ws.Range("I9").Value = 340.9
$ws.Range("K9").Value = 340.9
$ws.Range("M9").Value = -171.9

$ws.Range("H28").Value = 727.44446
$ws.Range("I28").Value = 658.86664
$ws.Range("K28").Value = 658.86664
$ws.Range("M28").Value = -173.86664

$ws.Range("H33").Value = 513.7083
$ws.Range("I33").Value = 150.6875
$ws.Range("K33").Value = 150.6875
$ws.Range("M33").Value = 78.3125

$ws.Range("H40").Value = 5083.1665
$ws.Range("J40").Value = 5083.1665
$ws.Range("L40").Value = 5083.1665
$ws.Range("N40").Value = -5433.1665

$ws.Range("H64").Value = 12105
$ws.Range("I64").Value = 9461.23
$ws.Range("J64").Value = 17833.166
$ws.Range("K64").Value = 9461.23
$ws.Range("L64").Value = 17833.166
$ws.Range("M64").Value = -9213.23
$ws.Range("N64").Value = -18329.166

$ws.Range("H67").Value = 12105
$ws.Range("I67").Value = 9461.23
$ws.Range("J67").Value = 17833.166
$ws.Range("K67").Value = 9461.23
$ws.Range("L67").Value = 17833.166
$ws.Range("M67").Value = -8603.23
$ws.Range("N67").Value = -19549.166

$ws.Range("H112").Value = 3263.2727
$ws.Range("I112").Value = 1248
$ws.Range("K112").Value = 3744
$ws.Range("M112").Value = -2636

$ws.Range("H116").Value = 6246.533
$ws.Range("I116").Value = 5659.7334
$ws.Range("K116").Value = 5659.7334
$ws.Range("M116").Value = -2217.7334

$ws.Range("H138").Value = 6824.05
$ws.Range("I138").Value = 4900.3125
$ws.Range("J138").Value = 8106.5415
$ws.Range("K138").Value = 14700.9375
$ws.Range("L138").Value = 24319.6245
$ws.Range("M138").Value = -9560.9375
$ws.Range("N138").Value = -34599.62450000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 420.375
$ws.Range("I2").Value = 420.375
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 420.375
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -307.375
$ws.Range("N2").ClearContents()

$ws.Range("H32").Value = 4600.7085
$ws.Range("I32").Value = 3496.5
$ws.Range("K32").Value = 3496.5
$ws.Range("M32").Value = -3209.5

$ws.Range("H45").Value = 2250
$ws.Range("I45").Value = 2333.3333
$ws.Range("K45").Value = 2333.3333
$ws.Range("M45").Value = -1956.3333

$ws.Range("H97").Value = 207.11111
$ws.Range("I97").Value = 207.11111
$ws.Range("K97").Value = 207.11111
$ws.Range("M97").Value = 288.88889

$ws.Range("H110").Value = 1820.4286
$ws.Range("I110").Value = 1554.2222
$ws.Range("J110").Value = 2299.6
$ws.Range("K110").Value = 1554.2222
$ws.Range("L110").Value = 2299.6
$ws.Range("M110").Value = 490.7778000000001
$ws.Range("N110").Value = -6389.6

$ws.Range("H116").Value = 420.375
$ws.Range("I116").Value = 420.375
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 420.375
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1873.625
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 3997.5386
$ws.Range("I132").Value = 2686.2856
$ws.Range("J132").Value = 5527.3335
$ws.Range("K132").Value = 8058.8568
$ws.Range("L132").Value = 16582.0005
$ws.Range("M132").Value = -5528.8568
$ws.Range("N132").Value = -21642.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 420.375
$ws.Range("I3").Value = 420.375
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 420.375
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -306.375
$ws.Range("N3").ClearContents()

$ws.Range("H76").Value = 20666.666
$ws.Range("J76").Value = 20666.666
$ws.Range("L76").Value = 20666.666
$ws.Range("N76").Value = -21296.666

$ws.Range("H79").Value = 20666.666
$ws.Range("J79").Value = 20666.666
$ws.Range("L79").Value = 20666.666
$ws.Range("N79").Value = -22850.666

$ws.Range("H105").Value = 4942
$ws.Range("I105").Value = 3198.1667
$ws.Range("K105").Value = 3198.1667
$ws.Range("M105").Value = -1451.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1319.2
$ws.Range("I16").Value = 1387.1428
$ws.Range("J16").Value = 1160.6666
$ws.Range("K16").Value = 1387.1428
$ws.Range("L16").Value = 1160.6666
$ws.Range("M16").Value = -1100.1428
$ws.Range("N16").Value = -1734.6666

$ws.Range("H41").Value = 12184.833
$ws.Range("I41").Value = 3277.5
$ws.Range("J41").Value = 29999.5
$ws.Range("K41").Value = 3277.5
$ws.Range("L41").Value = 29999.5
$ws.Range("M41").Value = -2849.5
$ws.Range("N41").Value = -30855.5

$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()

$ws.Range("H62").Value = 4003
$ws.Range("J62").Value = 3006
$ws.Range("L62").Value = 3006
$ws.Range("N62").Value = -4254

$ws.Range("H65").Value = 4003
$ws.Range("J65").Value = 3006
$ws.Range("L65").Value = 15030
$ws.Range("N65").Value = -21270

$ws.Range("H88").Value = 19648
$ws.Range("J88").Value = 19648
$ws.Range("L88").Value = 19648
$ws.Range("N88").Value = -20460

$ws.Range("H91").Value = 19648
$ws.Range("J91").Value = 19648
$ws.Range("L91").Value = 19648
$ws.Range("N91").Value = -22456

$ws.Range("H113").Value = 1319.2
$ws.Range("I113").Value = 1387.1428
$ws.Range("J113").Value = 1160.6666
$ws.Range("K113").Value = 1387.1428
$ws.Range("L113").Value = 1160.6666
$ws.Range("M113").Value = 782.8571999999999
$ws.Range("N113").Value = -5500.6666

$ws.Range("H134").Value = 3989.7
$ws.Range("I134").Value = 3766.4443
$ws.Range("K134").Value = 11299.3329
$ws.Range("M134").Value = -8764.332900000001

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 860.94446
$ws.Range("I34").Value = 170
$ws.Range("J34").Value = 1724.625
$ws.Range("K34").Value = 510
$ws.Range("L34").Value = 5173.875
$ws.Range("M34").Value = -426
$ws.Range("N34").Value = -5341.875

$ws.Range("H50").Value = 254.8
$ws.Range("I50").Value = 254.8
$ws.Range("K50").Value = 764.4000000000001
$ws.Range("M50").Value = -283.4000000000001

$ws.Range("H53").Value = 254.8
$ws.Range("I53").Value = 254.8
$ws.Range("K53").Value = 764.4000000000001
$ws.Range("M53").Value = -283.4000000000001

$ws.Range("H60").Value = 3743485
$ws.Range("I60").Value = 575
$ws.Range("J60").Value = 11229305
$ws.Range("K60").Value = 1725
$ws.Range("L60").Value = 33687915
$ws.Range("M60").Value = -1474
$ws.Range("N60").Value = -33688417

$ws.Range("H68").Value = 1999.75
$ws.Range("I68").Value = 1999.75
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 5999.25
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -5188.25
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1999.75
$ws.Range("I71").Value = 1999.75
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 17997.75
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -13941.75
$ws.Range("N71").ClearContents()

$ws.Range("H121").Value = 1369.25
$ws.Range("I121").Value = 825.6667
$ws.Range("K121").Value = 2477.0001
$ws.Range("M121").Value = -1167.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2746.923
$ws.Range("I80").Value = 2588.1667
$ws.Range("J80").Value = 2883
$ws.Range("K80").Value = 2588.1667
$ws.Range("L80").Value = 2883
$ws.Range("M80").Value = -1590.1667
$ws.Range("N80").Value = -4879

$ws.Range("H83").Value = 2746.923
$ws.Range("I83").Value = 2588.1667
$ws.Range("J83").Value = 2883
$ws.Range("K83").Value = 12940.8335
$ws.Range("L83").Value = 14415
$ws.Range("M83").Value = -7948.833500000001
$ws.Range("N83").Value = -24399

$ws.Range("H102").Value = 142857650
$ws.Range("I102").Value = 200000460
$ws.Range("J102").Value = 612.5
$ws.Range("K102").Value = 200000460
$ws.Range("L102").Value = 612.5
$ws.Range("M102").Value = -199998838
$ws.Range("N102").Value = -3856.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 60000
$ws.Range("I40").Value = 60000
$ws.Range("K40").Value = 60000
$ws.Range("M40").Value = -59851

$ws.Range("H80").Value = 30301
$ws.Range("J80").Value = 30301
$ws.Range("L80").Value = 30301
$ws.Range("N80").Value = -32297

$ws.Range("H83").Value = 30301
$ws.Range("J83").Value = 30301
$ws.Range("L83").Value = 90903
$ws.Range("N83").Value = -100887
